# Auto-update draw results: append the 2025-11-28 Pick 3 draw as a new
# row (row 73) at the bottom of the Results sheet, matching the existing
# data layout (dates/phase codes stored as literal text, not numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.Range("A73:E73")

# Format the new cells as Text *before* assigning values so values that
# look numeric/date-like (e.g. "2025-11-28", "251128") are stored as
# literal text, consistent with every other row in the sheet.
$newRow.NumberFormat = "@"

$ws.Range("A73").Value = "2025-11-28"
$ws.Range("B73").Value = "Pick 3"
$ws.Range("C73").Value = "251128"
$ws.Range("D73").Value = "3-0-3"
$ws.Range("E73").Value = "2025-11-28T21:38:19.934+04:00"
